# Update the "Urbanization" breakdown labels on the active sheet.
# Old rows 23/24 ("Urban"/"Rural" trio) are replaced with the more
# specific "City"/"Village" trio (Kyrgyz / Russian / English columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23: Urban -> City
$ws.Range("A23").Value = "Шаар жерлери"
$ws.Range("B23").Value = "Городские поселения"
$ws.Range("C23").Value = "City"

# Row 24: Rural -> Village
$ws.Range("A24").Value = "Айыл аймагы"
$ws.Range("B24").Value = "Сельская местность"
$ws.Range("C24").Value = "Village"

# Move the active selection to C30, matching the saved view state.
$ws.Range("C30").Select()
